$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 502.55554
$ws.Cells.Item(2,10).Value = 689
$ws.Cells.Item(2,12).Value = 689
$ws.Cells.Item(2,14).Value = -915
$ws.Cells.Item(11,8).Value = 33.875
$ws.Cells.Item(11,9).Value = 33.875
$ws.Cells.Item(11,11).Value = 33.875
$ws.Cells.Item(11,13).Value = 106.125
$ws.Cells.Item(18,8).Value = 0
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,12).ClearContents()
$ws.Cells.Item(18,14).Value = 0
$ws.Cells.Item(31,8).Value = 90.818184
$ws.Cells.Item(31,9).Value = 69.5
$ws.Cells.Item(31,11).Value = 208.5
$ws.Cells.Item(31,13).Value = 21.5
$ws.Cells.Item(43,8).Value = 11749
$ws.Cells.Item(43,9).Value = 8998
$ws.Cells.Item(43,11).Value = 8998
$ws.Cells.Item(43,13).Value = -8929
$ws.Cells.Item(62,8).Value = 2900.8333
$ws.Cells.Item(62,9).Value = 1851.5
$ws.Cells.Item(62,10).Value = 4999.5
$ws.Cells.Item(62,11).Value = 1851.5
$ws.Cells.Item(62,12).Value = 4999.5
$ws.Cells.Item(62,13).Value = -1227.5
$ws.Cells.Item(62,14).Value = -6247.5
$ws.Cells.Item(64,8).Value = 2666
$ws.Cells.Item(65,8).Value = 2900.8333
$ws.Cells.Item(65,9).Value = 1851.5
$ws.Cells.Item(65,10).Value = 4999.5
$ws.Cells.Item(65,11).Value = 9257.5
$ws.Cells.Item(65,12).Value = 24997.5
$ws.Cells.Item(65,13).Value = -6137.5
$ws.Cells.Item(65,14).Value = -31237.5
$ws.Cells.Item(67,8).Value = 2666
$ws.Cells.Item(94,8).Value = 3131.6924
$ws.Cells.Item(94,9).Value = 3267.6667
$ws.Cells.Item(94,10).Value = 1500
$ws.Cells.Item(94,11).Value = 3267.6667
$ws.Cells.Item(94,12).Value = 1500
$ws.Cells.Item(94,13).Value = -2816.6667
$ws.Cells.Item(94,14).Value = -2402
$ws.Cells.Item(98,8).Value = 1279.8823
$ws.Cells.Item(98,9).Value = 1254.5
$ws.Cells.Item(98,11).Value = 1254.5
$ws.Cells.Item(98,13).Value = 243.5
$ws.Cells.Item(116,8).Value = 3723.75
$ws.Cells.Item(116,9).Value = 2965
$ws.Cells.Item(116,11).Value = 2965
$ws.Cells.Item(116,13).Value = 477
$ws.Cells.Item(122,8).Value = 1279.8823
$ws.Cells.Item(122,9).Value = 1254.5
$ws.Cells.Item(122,11).Value = 3763.5
$ws.Cells.Item(122,13).Value = -1313.5
$ws.Cells.Item(141,8).Value = 4398
$ws.Cells.Item(141,9).Value = 3997.5
$ws.Cells.Item(141,11).Value = 11992.5
$ws.Cells.Item(141,13).Value = -6812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10,8).Value = 12750
$ws.Cells.Item(10,10).Value = 5500
$ws.Cells.Item(10,12).Value = 5500
$ws.Cells.Item(10,14).Value = -5840
$ws.Cells.Item(21,8).Value = 0
$ws.Cells.Item(21,10).Value = 0
$ws.Cells.Item(21,12).ClearContents()
$ws.Cells.Item(21,14).Value = 0
$ws.Cells.Item(32,8).Value = 3866.973
$ws.Cells.Item(32,9).Value = 2517.5757
$ws.Cells.Item(32,11).Value = 2517.5757
$ws.Cells.Item(32,13).Value = -2230.5757
$ws.Cells.Item(33,8).Value = 0
$ws.Cells.Item(33,9).Value = 0
$ws.Cells.Item(33,11).Value = 0
$ws.Cells.Item(33,13).ClearContents()
$ws.Cells.Item(36,8).Value = 2466.6667
$ws.Cells.Item(36,9).Value = 2466.6667
$ws.Cells.Item(36,11).Value = 2466.6667
$ws.Cells.Item(36,13).Value = -2120.6667
$ws.Cells.Item(61,8).Value = 6233.7144
$ws.Cells.Item(61,9).Value = 6328.615
$ws.Cells.Item(61,11).Value = 6328.615
$ws.Cells.Item(61,13).Value = -6116.615
$ws.Cells.Item(74,8).Value = 1377.7
$ws.Cells.Item(74,9).Value = 893.4
$ws.Cells.Item(74,11).Value = 893.4
$ws.Cells.Item(74,13).Value = -19.39999999999998
$ws.Cells.Item(77,8).Value = 1377.7
$ws.Cells.Item(77,9).Value = 893.4
$ws.Cells.Item(77,11).Value = 4467
$ws.Cells.Item(77,13).Value = -99
$ws.Cells.Item(136,8).Value = 6233.7144
$ws.Cells.Item(136,9).Value = 6328.615
$ws.Cells.Item(136,11).Value = 18985.845
$ws.Cells.Item(136,13).Value = -16435.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value = 1056.8572
$ws.Cells.Item(86,9).Value = 849.5
$ws.Cells.Item(86,10).Value = 1333.3334
$ws.Cells.Item(86,11).Value = 849.5
$ws.Cells.Item(86,12).Value = 1333.3334
$ws.Cells.Item(86,13).Value = 273.5
$ws.Cells.Item(86,14).Value = -3579.3334
$ws.Cells.Item(89,8).Value = 1056.8572
$ws.Cells.Item(89,9).Value = 849.5
$ws.Cells.Item(89,10).Value = 1333.3334
$ws.Cells.Item(89,11).Value = 4247.5
$ws.Cells.Item(89,12).Value = 6666.666999999999
$ws.Cells.Item(89,13).Value = 1368.5
$ws.Cells.Item(89,14).Value = -17898.667
$ws.Cells.Item(107,8).Value = 1909
$ws.Cells.Item(107,9).Value = 2133
$ws.Cells.Item(107,10).Value = 1013
$ws.Cells.Item(107,11).Value = 2133
$ws.Cells.Item(107,12).Value = 1013
$ws.Cells.Item(107,13).Value = -213
$ws.Cells.Item(107,14).Value = -4853

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7,8).Value = 124.05556
$ws.Cells.Item(7,9).Value = 102.75
$ws.Cells.Item(7,11).Value = 102.75
$ws.Cells.Item(7,13).Value = 10.25
$ws.Cells.Item(22,8).Value = 101796.08
$ws.Cells.Item(22,9).Value = 168549.72
$ws.Cells.Item(22,11).Value = 168549.72
$ws.Cells.Item(22,13).Value = -168199.72
$ws.Cells.Item(74,8).Value = 0
$ws.Cells.Item(74,10).Value = 0
$ws.Cells.Item(74,12).ClearContents()
$ws.Cells.Item(74,14).Value = 0
$ws.Cells.Item(77,8).Value = 0
$ws.Cells.Item(77,10).Value = 0
$ws.Cells.Item(77,12).ClearContents()
$ws.Cells.Item(77,14).Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(82,8).Value = 0
$ws.Cells.Item(82,9).Value = 0
$ws.Cells.Item(82,11).Value = 0
$ws.Cells.Item(82,13).ClearContents()
$ws.Cells.Item(85,8).Value = 0
$ws.Cells.Item(85,9).Value = 0
$ws.Cells.Item(85,11).Value = 0
$ws.Cells.Item(85,13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14,8).Value = 100
$ws.Cells.Item(14,9).Value = 100
$ws.Cells.Item(14,10).Value = 0
$ws.Cells.Item(14,11).Value = 100
$ws.Cells.Item(14,12).ClearContents()
$ws.Cells.Item(14,13).Value = 68
$ws.Cells.Item(14,14).Value = 0
$ws.Cells.Item(20,8).Value = 23000
$ws.Cells.Item(20,9).Value = 0
$ws.Cells.Item(20,10).Value = 23000
$ws.Cells.Item(20,11).Value = 0
$ws.Cells.Item(20,12).ClearContents()
$ws.Cells.Item(20,13).Value = 23000
$ws.Cells.Item(20,14).Value = -23490
$ws.Cells.Item(24,8).Value = 16500
$ws.Cells.Item(24,9).Value = 0
$ws.Cells.Item(24,10).Value = 16500
$ws.Cells.Item(24,11).Value = 0
$ws.Cells.Item(24,12).ClearContents()
$ws.Cells.Item(24,13).Value = 16500
$ws.Cells.Item(24,14).Value = -16846
$ws.Cells.Item(70,8).Value = 24997.5
$ws.Cells.Item(70,9).Value = 24997
$ws.Cells.Item(70,11).Value = 24997
$ws.Cells.Item(70,13).Value = -24727
$ws.Cells.Item(73,8).Value = 24997.5
$ws.Cells.Item(73,9).Value = 24997
$ws.Cells.Item(73,11).Value = 24997
$ws.Cells.Item(73,13).Value = -24061
$ws.Cells.Item(80,8).Value = 4500
$ws.Cells.Item(80,9).Value = 0
$ws.Cells.Item(80,10).Value = 4500
$ws.Cells.Item(80,11).Value = 0
$ws.Cells.Item(80,12).ClearContents()
$ws.Cells.Item(80,13).Value = 4500
$ws.Cells.Item(80,14).Value = -6496
$ws.Cells.Item(83,8).Value = 4500
$ws.Cells.Item(83,9).Value = 0
$ws.Cells.Item(83,10).Value = 4500
$ws.Cells.Item(83,11).Value = 0
$ws.Cells.Item(83,12).ClearContents()
$ws.Cells.Item(83,13).Value = 22500
$ws.Cells.Item(83,14).Value = -32484
$ws.Cells.Item(105,8).Value = 16166.667
$ws.Cells.Item(105,10).Value = 16166.667
$ws.Cells.Item(105,12).Value = 16166.667
$ws.Cells.Item(105,14).Value = -23154.667
$ws.Cells.Item(132,8).Value = 2317.439
$ws.Cells.Item(132,9).Value = 2430.9722
$ws.Cells.Item(132,10).Value = 1500
$ws.Cells.Item(132,11).Value = 7292.9166
$ws.Cells.Item(132,12).Value = 4500
$ws.Cells.Item(132,13).Value = -4762.9166
$ws.Cells.Item(132,14).Value = -9560

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17,8).Value = 18504.5
$ws.Cells.Item(17,9).Value = 18000
$ws.Cells.Item(17,10).Value = 19009
$ws.Cells.Item(17,11).Value = 18000
$ws.Cells.Item(17,12).Value = 19009
$ws.Cells.Item(17,13).Value = -17830
$ws.Cells.Item(17,14).Value = -19349
$ws.Cells.Item(18,8).Value = 59999.5
$ws.Cells.Item(18,10).Value = 59999.5
$ws.Cells.Item(18,12).Value = 59999.5
$ws.Cells.Item(18,14).Value = -60343.5
$ws.Cells.Item(22,8).Value = 2619.5217
$ws.Cells.Item(22,10).Value = 3885
$ws.Cells.Item(22,12).Value = 3885
$ws.Cells.Item(22,14).Value = -4475
$ws.Cells.Item(27,8).Value = 2619.5217
$ws.Cells.Item(27,10).Value = 3885
$ws.Cells.Item(27,12).Value = 3885
$ws.Cells.Item(27,14).Value = -4099

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12,8).Value = 15413.667
$ws.Cells.Item(12,9).Value = 15997.5
$ws.Cells.Item(12,10).Value = 14246
$ws.Cells.Item(12,11).Value = 15997.5
$ws.Cells.Item(12,12).Value = 14246
$ws.Cells.Item(12,13).Value = -15855.5
$ws.Cells.Item(12,14).Value = -14530
$ws.Cells.Item(17,8).Value = 33333884
$ws.Cells.Item(17,9).Value = 50000224
$ws.Cells.Item(17,11).Value = 50000224
$ws.Cells.Item(17,13).Value = -50000052
$ws.Cells.Item(20,9).Value = 0
$ws.Cells.Item(20,11).Value = 0
$ws.Cells.Item(20,13).ClearContents()
$ws.Cells.Item(31,8).Value = 21000.6
$ws.Cells.Item(31,10).Value = 21000.6
$ws.Cells.Item(31,12).Value = 21000.6
$ws.Cells.Item(31,14).Value = -21696.6
$ws.Cells.Item(51,8).Value = 0
$ws.Cells.Item(51,9).Value = 0
$ws.Cells.Item(51,11).Value = 0
$ws.Cells.Item(51,13).ClearContents()
$ws.Cells.Item(52,8).Value = 0
$ws.Cells.Item(52,9).Value = 0
$ws.Cells.Item(52,11).Value = 0
$ws.Cells.Item(52,13).ClearContents()
$ws.Cells.Item(122,8).Value = 2301.926
$ws.Cells.Item(122,9).Value = 1845.5294
$ws.Cells.Item(122,10).Value = 3077.8
$ws.Cells.Item(122,11).Value = 5536.5882
$ws.Cells.Item(122,12).Value = 9233.400000000001
$ws.Cells.Item(122,13).Value = -3086.5882
$ws.Cells.Item(122,14).Value = -14133.4
